$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 22.14384766666667
$ws.Range("H2").Value = 66.431543
$ws.Range("I2").Value = 0.05562336639723622
$ws.Range("J2").Value = 0.0556233663972362
$ws.Range("O2").Value = 0.3099390012751145
$ws.Range("P2").Value = 0.3099390012751145
$ws.Range("Q2").Value = 2.655342586535556
$ws.Range("R2").Value = 23.89808327882
$ws.Range("S2").Value = 0.01723985062871916
$ws.Range("T2").Value = 0.01723985062871915

$ws.Range("G3").Value = 22.14384766666667
$ws.Range("H3").Value = 66.431543
$ws.Range("I3").Value = 0.05562336639723622
$ws.Range("J3").Value = 0.0556233663972362
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2669800000000001
$ws.Range("N3").Value = 0.8009400000000001
$ws.Range("O3").Value = 0.6900609987248855
$ws.Range("P3").Value = 0.6900609987248854
$ws.Range("Q3").Value = 5.911964450046669
$ws.Range("R3").Value = 53.20768005042001
$ws.Range("S3").Value = 0.03838351576851706
$ws.Range("T3").Value = 0.03838351576851704

$ws.Range("I4").Value = 0.8709978578802913
$ws.Range("J4").Value = 0.8709978578802913
$ws.Range("O4").Value = 0.3099390012751145
$ws.Range("P4").Value = 0.3099390012751145
$ws.Range("S4").Value = 0.2699562061841816
$ws.Range("T4").Value = 0.2699562061841816

$ws.Range("I5").Value = 0.8709978578802913
$ws.Range("J5").Value = 0.8709978578802913
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2669800000000001
$ws.Range("N5").Value = 0.8009400000000001
$ws.Range("O5").Value = 0.6900609987248855
$ws.Range("P5").Value = 0.6900609987248854
$ws.Range("Q5").Value = 92.57455464096002
$ws.Range("R5").Value = 833.1709917686401
$ws.Range("S5").Value = 0.6010416516961097
$ws.Range("T5").Value = 0.6010416516961096

$ws.Range("G6").Value = 0.1541363333333333
$ws.Range("H6").Value = 0.462409
$ws.Range("I6").Value = 0.0003871766945467397
$ws.Range("J6").Value = 0.0003871766945467395
$ws.Range("O6").Value = 0.3099390012751145
$ws.Range("P6").Value = 0.3099390012751145
$ws.Range("Q6").Value = 0.01848300151777778
$ws.Range("R6").Value = 0.16634701366
$ws.Range("S6").Value = 0.0001200011580248166
$ws.Range("T6").Value = 0.0001200011580248165

$ws.Range("G7").Value = 0.1541363333333333
$ws.Range("H7").Value = 0.462409
$ws.Range("I7").Value = 0.0003871766945467397
$ws.Range("J7").Value = 0.0003871766945467395
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2669800000000001
$ws.Range("N7").Value = 0.8009400000000001
$ws.Range("O7").Value = 0.6900609987248855
$ws.Range("P7").Value = 0.6900609987248854
$ws.Range("Q7").Value = 0.04115131827333335
$ws.Range("R7").Value = 0.3703618644600001
$ws.Range("S7").Value = 0.0002671755365219231
$ws.Range("T7").Value = 0.0002671755365219229

$ws.Range("G8").Value = 28.90575466666667
$ws.Range("H8").Value = 86.717264
$ws.Range("I8").Value = 0.07260867248616912
$ws.Range("J8").Value = 0.07260867248616912
$ws.Range("O8").Value = 0.3099390012751145
$ws.Range("P8").Value = 0.3099390012751145
$ws.Range("Q8").Value = 3.466185394595556
$ws.Range("R8").Value = 31.19566855136
$ws.Range("S8").Value = 0.02250425943427514
$ws.Range("T8").Value = 0.02250425943427514

$ws.Range("G9").Value = 28.90575466666667
$ws.Range("H9").Value = 86.717264
$ws.Range("I9").Value = 0.07260867248616912
$ws.Range("J9").Value = 0.07260867248616912
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2669800000000001
$ws.Range("N9").Value = 0.8009400000000001
$ws.Range("O9").Value = 0.6900609987248855
$ws.Range("P9").Value = 0.6900609987248854
$ws.Range("Q9").Value = 7.717258380906668
$ws.Range("R9").Value = 69.45532542816001
$ws.Range("S9").Value = 0.05010441305189398
$ws.Range("T9").Value = 0.05010441305189397

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1524443333333333
$ws.Range("H10").Value = 0.457333
$ws.Range("I10").Value = 0.0003829265417566354
$ws.Range("J10").Value = 0.0003829265417566354
$ws.Range("O10").Value = 0.3099390012751145
$ws.Range("P10").Value = 0.3099390012751145
$ws.Range("Q10").Value = 0.01828010815777778
$ws.Range("R10").Value = 0.16452097342
$ws.Range("S10").Value = 0.000118683869913785
$ws.Range("T10").Value = 0.000118683869913785

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1524443333333333
$ws.Range("H11").Value = 0.457333
$ws.Range("I11").Value = 0.0003829265417566354
$ws.Range("J11").Value = 0.0003829265417566354
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2669800000000001
$ws.Range("N11").Value = 0.8009400000000001
$ws.Range("O11").Value = 0.6900609987248855
$ws.Range("P11").Value = 0.6900609987248854
$ws.Range("Q11").Value = 0.04069958811333334
$ws.Range("R11").Value = 0.3662962930200001
$ws.Range("S11").Value = 0.0002642426718428504
$ws.Range("T11").Value = 0.0002642426718428504
